$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column H first (sectionname column), then the F/G pairs row by row -
# matches the order in which the new shared strings were introduced by
# the original author.
$ws.Range("H2").Value = "pop_filter1_section"
$ws.Range("H3").Value = "pop_filter1_section"
$ws.Range("H4").Value = "pop_filter2_section"

$ws.Range("F2").Value = "pop_filter1_section1"
$ws.Range("G2").Value = "pop_filter1_section1_checkbox"

$ws.Range("F3").Value = "pop_filter1_section2"
$ws.Range("G3").Value = "pop_filter1_section2_checkbox"

$ws.Range("F4").Value = "pop_filter2_section2"
$ws.Range("G4").Value = "pop_filter2_section2_checkbox"

$ws.Activate()
$ws.Range("F5").Select()
$excel.ActiveWindow.ScrollColumn = 3
